$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top, pushing existing data down to row 2
$ws.Rows.Item(1).Insert()

# Populate the new header row with titles
$ws.Range("A1").Value = "Type of lecture"
$ws.Range("B1").Value = "Lecture Title"
$ws.Range("C1").Value = "Resource Person Name"
$ws.Range("D1").Value = "Designation"
$ws.Range("E1").Value = "Company Name"
$ws.Range("F1").Value = "City,Country"
$ws.Range("G1").Value = "Course Code/Course Title"
$ws.Range("H1").Value = "Date"
$ws.Range("I1").Value = "Time"
$ws.Range("J1").Value = "Mode/Venue"
$ws.Range("K1").Value = "Coordinators"
